$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")

# Use higher value for plant types early in learning curve
$ws.Range("B19:B25").Value = 0.75

# Update the selection to reflect the rows that were just edited
$ws.Activate()
$ws.Range("B19:B25").Select()
